$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for the columns that get shuffled
# across the data rows (2..10): D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg).
$orig = @{}
for ($r = 2; $r -le 10; $r++) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

# Map each target row to the source row whose original values it should
# receive (derived from the diff between before/after states).
$mapping = @{
    2  = 6
    3  = 10
    4  = 2
    5  = 7
    6  = 3
    7  = 8
    8  = 9
    9  = 5
    10 = 4
}

foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    $src = $orig[$srcRow]

    $ws.Cells.Item($targetRow, 4).Value2 = $src.D
    $ws.Cells.Item($targetRow, 13).Value2 = $src.M
    $ws.Cells.Item($targetRow, 14).Value2 = $src.N
    $ws.Cells.Item($targetRow, 15).Value2 = $src.O
    $ws.Cells.Item($targetRow, 16).Value2 = $src.P
    $ws.Cells.Item($targetRow, 19).Value2 = $src.S
}

$wb.Save()
